$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 254 (shifts existing rows 254:288 down to 256:290)
$ws.Rows.Item(254).Insert()
$ws.Rows.Item(254).Insert()

# New row 254
$ws.Cells.Item(254, 1).Value = 7
$ws.Cells.Item(254, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(254, 3).Value = "Ñuble"
$ws.Cells.Item(254, 4).Value = 44984
$ws.Cells.Item(254, 5).Value = 16
$ws.Cells.Item(254, 6).Value = 100112024
$ws.Cells.Item(254, 7).Value = "Choclo"
$ws.Cells.Item(254, 8).Value = "Choclero"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 30000
$ws.Cells.Item(254, 11).Value = 400
$ws.Cells.Item(254, 12).Value = 450
$ws.Cells.Item(254, 13).Value = 425
$ws.Cells.Item(254, 14).Value = "$/unidad"
$ws.Cells.Item(254, 15).Value = "Región del Maule"
$ws.Cells.Item(254, 16).Value = 425
$ws.Cells.Item(254, 17).Value = 1
$ws.Cells.Item(254, 18).Value = "Hortaliza"

# New row 255
$ws.Cells.Item(255, 1).Value = 7
$ws.Cells.Item(255, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(255, 3).Value = "Ñuble"
$ws.Cells.Item(255, 4).Value = 44984
$ws.Cells.Item(255, 5).Value = 16
$ws.Cells.Item(255, 6).Value = 100112024
$ws.Cells.Item(255, 7).Value = "Choclo"
$ws.Cells.Item(255, 8).Value = "Choclero"
$ws.Cells.Item(255, 9).Value = "Segunda"
$ws.Cells.Item(255, 10).Value = 15000
$ws.Cells.Item(255, 11).Value = 350
$ws.Cells.Item(255, 12).Value = 350
$ws.Cells.Item(255, 13).Value = 350
$ws.Cells.Item(255, 14).Value = "$/unidad"
$ws.Cells.Item(255, 15).Value = "Región del Maule"
$ws.Cells.Item(255, 16).Value = 350
$ws.Cells.Item(255, 17).Value = 1
$ws.Cells.Item(255, 18).Value = "Hortaliza"
